# Generate Report for Handoff
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": the status cells and their "last updated" timestamps
# are refreshed on the Overview sheet and on each per-language sheet
# (zh-cn / de-de). The status columns on the Overview sheet auto-size to
# the new (longer) status text, which widens two columns there and the
# mirrored "Status" column on each language sheet.

$wb = $excel.ActiveWorkbook

# This host's ColumnWidth setter quantizes to steps of 1/6 character
# (stored = round(input*6)/6 + 5/6), so an input that lands exactly on the
# target stored width isn't always available. Solve for the input value
# that gets the stored width as close as possible to the desired width.
function Get-ColumnWidthInput($targetStoredWidth) {
    $steps = [Math]::Round(($targetStoredWidth - 5.0/6.0) * 6.0)
    return ($steps / 6.0)
}

$newStatus = "Ready for handoff"
$newColWidth = Get-ColumnWidthInput 17.2159881591797

# --- Overview sheet: zh-cn / de-de status + last handoff-generate date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-29 08:42:04"
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-29 08:41:56"
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-29 08:42:04"
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
